# Lab2/Histograms.xlsx edit:
#  - rename sheets img1 -> img0, img2 -> img1
#  - delete the (empty) Plan3 sheet
#  - move the "active sheet" from img2 (was index 1 / img2) to img1 (new index 0 / img0)
#  - update both histogram chart titles / series / source refs to match the new sheet names
#  - tidy up the sheet view selections that Excel re-records when the workbook is saved

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- rename worksheets (order matters: rename img1 -> img0 first so the name
#     "img1" is free again before img2 takes it) -------------------------------
$wb.Worksheets("img1").Name = "img0"
$wb.Worksheets("img2").Name = "img1"

# --- drop the unused third sheet ----------------------------------------------
$wb.Worksheets("Plan3").Delete()

$wsImg0 = $wb.Worksheets("img0")
$wsImg1 = $wb.Worksheets("img1")

# --- chart on img0 (was img1 Histogram, now image0 Histogram) ----------------
$chart0 = $wsImg0.ChartObjects(1).Chart
$chart0.HasTitle = $true
$chart0.ChartTitle.Text = "image0 Histogram"
$ser0 = $chart0.SeriesCollection(1)
$ser0.Formula = '=SERIES("image0 Histogram",img0!$A$2:$A$257,img0!$B$2:$B$257,1)'

# --- chart on img1 (was img2 Histogram, now image1 Histogram) ----------------
$chart1 = $wsImg1.ChartObjects(1).Chart
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "image1 Histogram"
$ser1 = $chart1.SeriesCollection(1)
$ser1.Formula = "=SERIES(""image1 Histogram"",'img1'!`$A`$2:`$A`$257,'img1'!`$B`$2:`$B`$257,1)"

# --- sheet views: img0 becomes the active / selected sheet --------------------
$wsImg0.Activate()
$wsImg0.Range("Q20").Select()

$wsImg1.Range("M24").Select()
